$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.01"
$ws.Range("D2").NumberFormat = "general"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.66%"
$ws.Range("E2").NumberFormat = "general"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.60"
$ws.Range("D3").NumberFormat = "general"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.79%"
$ws.Range("E3").NumberFormat = "general"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.046"
$ws.Range("D4").NumberFormat = "general"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.64%"
$ws.Range("E4").NumberFormat = "general"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07976"
$ws.Range("D5").NumberFormat = "general"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.55%"
$ws.Range("E5").NumberFormat = "general"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.49%"
$ws.Range("E6").NumberFormat = "general"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("B7").NumberFormat = "general"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C7").NumberFormat = "general"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.775"
$ws.Range("D7").NumberFormat = "general"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.07%"
$ws.Range("E7").NumberFormat = "general"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("B8").NumberFormat = "general"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C8").NumberFormat = "general"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9238"
$ws.Range("D8").NumberFormat = "general"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.72%"
$ws.Range("E8").NumberFormat = "general"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B9").NumberFormat = "general"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C9").NumberFormat = "general"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1282"
$ws.Range("D9").NumberFormat = "general"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.82%"
$ws.Range("E9").NumberFormat = "general"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("B10").NumberFormat = "general"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").NumberFormat = "general"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1885"
$ws.Range("D10").NumberFormat = "general"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.04%"
$ws.Range("E10").NumberFormat = "general"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("B11").NumberFormat = "general"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").NumberFormat = "general"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09023"
$ws.Range("D11").NumberFormat = "general"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.70%"
$ws.Range("E11").NumberFormat = "general"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("B12").NumberFormat = "general"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").NumberFormat = "general"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03423"
$ws.Range("D12").NumberFormat = "general"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.16%"
$ws.Range("E12").NumberFormat = "general"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("B13").NumberFormat = "general"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").NumberFormat = "general"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09865"
$ws.Range("D13").NumberFormat = "general"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.04%"
$ws.Range("E13").NumberFormat = "general"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("B14").NumberFormat = "general"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").NumberFormat = "general"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001404"
$ws.Range("D14").NumberFormat = "general"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.98%"
$ws.Range("E14").NumberFormat = "general"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("B15").NumberFormat = "general"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").NumberFormat = "general"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006162"
$ws.Range("D15").NumberFormat = "general"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "6.55%"
$ws.Range("E15").NumberFormat = "general"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("B16").NumberFormat = "general"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C16").NumberFormat = "general"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.844"
$ws.Range("D16").NumberFormat = "general"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "7.85%"
$ws.Range("E16").NumberFormat = "general"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("B17").NumberFormat = "general"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C17").NumberFormat = "general"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.125"
$ws.Range("D17").NumberFormat = "general"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.86%"
$ws.Range("E17").NumberFormat = "general"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.377"
$ws.Range("D18").NumberFormat = "general"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "14.40%"
$ws.Range("E18").NumberFormat = "general"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3413"
$ws.Range("D19").NumberFormat = "general"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.93%"
$ws.Range("E19").NumberFormat = "general"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1311"
$ws.Range("D20").NumberFormat = "general"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.81%"
$ws.Range("E20").NumberFormat = "general"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.793"
$ws.Range("D21").NumberFormat = "general"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.09%"
$ws.Range("E21").NumberFormat = "general"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2309"
$ws.Range("D22").NumberFormat = "general"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-8.21%"
$ws.Range("E22").NumberFormat = "general"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04432"
$ws.Range("D23").NumberFormat = "general"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.59%"
$ws.Range("E23").NumberFormat = "general"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("D24").NumberFormat = "general"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.89%"
$ws.Range("E24").NumberFormat = "general"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004859"
$ws.Range("D25").NumberFormat = "general"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.80%"
$ws.Range("E25").NumberFormat = "general"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001302"
$ws.Range("D27").NumberFormat = "general"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-21.12%"
$ws.Range("E27").NumberFormat = "general"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "42.14%"
$ws.Range("E28").NumberFormat = "general"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01925"
$ws.Range("D39").NumberFormat = "general"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.55%"
$ws.Range("E39").NumberFormat = "general"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05146"
$ws.Range("D40").NumberFormat = "general"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.57%"
$ws.Range("E40").NumberFormat = "general"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007598"
$ws.Range("D41").NumberFormat = "general"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.49%"
$ws.Range("E41").NumberFormat = "general"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.52%"
$ws.Range("E42").NumberFormat = "general"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.23%"
$ws.Range("E43").NumberFormat = "general"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002133"
$ws.Range("D44").NumberFormat = "general"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.54%"
$ws.Range("E44").NumberFormat = "general"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009875"
$ws.Range("D45").NumberFormat = "general"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.21%"
$ws.Range("E45").NumberFormat = "general"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006213"
$ws.Range("D46").NumberFormat = "general"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.62%"
$ws.Range("E46").NumberFormat = "general"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"
$ws.Range("E47").NumberFormat = "general"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.01"
$ws.Range("D48").NumberFormat = "general"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.32%"
$ws.Range("E48").NumberFormat = "general"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.45%"
$ws.Range("E49").NumberFormat = "general"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D50").NumberFormat = "general"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
$ws.Range("E50").NumberFormat = "general"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("D51").NumberFormat = "general"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"
$ws.Range("E51").NumberFormat = "general"
